# Apply the edits described in the diff:
#  1. Update the invoice date field from 17/01/2014 to 25/01/2014
#     (the date lives inside a DOCPROPERTY fldSimple field, so update the
#      field's result text rather than Find/Replace, which does not
#      descend into field results)
#  2. Update the product description line 1 from "Dibond 2mm vari colori"
#     to "Supporti rigidi Forex 5mm bianco"
#  3. Update the product description line 2 from "Stampa UV alta qualità"
#     to "Stampa UV bassa qualità"
#  4. Update the quantity from "100" to "3"
#  5. Update the unit price from "1.718,6500" to "10,6000"
#  6. Update the total price from "171.865,00" to "31,80"

$d = $word.ActiveDocument

foreach ($fld in $d.Fields) {
    if ($fld.Code.Text -like "*invoice_date*") {
        $fld.Result.Text = "25/01/2014"
    }
}

$d.Content.Find.Execute("Dibond 2mm vari colori", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Supporti rigidi Forex 5mm bianco", 2)

$d.Content.Find.Execute("Stampa UV alta qualità", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Stampa UV bassa qualità", 2)

$d.Content.Find.Execute("100", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3", 2)

$d.Content.Find.Execute("1.718,6500", $true, $false, $false, $false, $false,
                         $true, 1, $false, "10,6000", 2)

$d.Content.Find.Execute("171.865,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "31,80", 2)
